$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

$ws.Range("A3").Value = "LOG020"
$ws.Range("B3").Value = "log allotment"

$ws.Activate()
$ws.Range("A4").Select()
